# Auto-generated edit script: updates HotStock_Top20 rankings (rows 2-21, cols A-C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "天赐材料"
$ws.Range("B2").Value = "中芯国际"

$ws.Range("A3").Value = "国轩高科"
$ws.Range("B3").Value = "国轩高科"
$ws.Range("C3").Value = "国轩高科"

$ws.Range("A4").Value = "三维通信"
$ws.Range("B4").Value = "三维通信"
$ws.Range("C4").Value = "吉视传媒"

$ws.Range("A5").Value = "多氟多"
$ws.Range("B5").Value = "天赐材料"
$ws.Range("C5").Value = "三维通信"

$ws.Range("A6").Value = "中芯国际"
$ws.Range("B6").Value = "卧龙电驱"
$ws.Range("C6").Value = "天赐材料"

$ws.Range("A7").Value = "卧龙电驱"
$ws.Range("B7").Value = "多氟多"
$ws.Range("C7").Value = "幸福蓝海"

$ws.Range("A8").Value = "先导智能"
$ws.Range("B8").Value = "天际股份"

$ws.Range("A9").Value = "天际股份"
$ws.Range("B9").Value = "三花智控"
$ws.Range("C9").Value = "万通发展"

$ws.Range("A10").Value = "三花智控"
$ws.Range("B10").Value = "晓程科技"
$ws.Range("C10").Value = "胜宏科技"

$ws.Range("A11").Value = "吉视传媒"
$ws.Range("B11").Value = "中国卫星"
$ws.Range("C11").Value = "步步高"

$ws.Range("A12").Value = "中国卫星"
$ws.Range("B12").Value = "吉视传媒"
$ws.Range("C12").Value = "工业富联"

$ws.Range("A13").Value = "步步高"
$ws.Range("B13").Value = "江特电机"
$ws.Range("C13").Value = "岩山科技"

$ws.Range("A14").Value = "万通发展"
$ws.Range("B14").Value = "春兴精工"
$ws.Range("C14").Value = "多氟多"

$ws.Range("A15").Value = "西部黄金"
$ws.Range("B15").Value = "西部黄金"
$ws.Range("C15").Value = "阳光电源"

$ws.Range("A16").Value = "岩山科技"
$ws.Range("B16").Value = "寒武纪-U"
$ws.Range("C16").Value = "济民健康"

$ws.Range("A17").Value = "晓程科技"
$ws.Range("B17").Value = "先导智能"
$ws.Range("C17").Value = "粤传媒"

$ws.Range("A18").Value = "上海洗霸"
$ws.Range("C18").Value = "金发科技"

$ws.Range("A19").Value = "江特电机"
$ws.Range("B19").Value = "步步高"
$ws.Range("C19").Value = "中芯国际"

$ws.Range("A20").Value = "四维图新"
$ws.Range("B20").Value = "首开股份"
$ws.Range("C20").Value = "中国电影"

$ws.Range("A21").Value = "阳光电源"
$ws.Range("B21").Value = "供销大集"
$ws.Range("C21").Value = "露笑科技"
